# Update the date line (unique text -> safe to use Find/Replace).
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-05-31 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-06-01 Saturday", 2) | Out-Null

# Update the 25 division-problem cells in the single table by explicit
# (row, column) address, so that the several values that collide with
# each other during the update (e.g. "58÷2=29, 0" moves from cell (1,5)
# to (5,3) while a different value moves into (1,2)) are never ambiguous.
$t = $d.Tables.Item(1)

$newValues = @{
    "1,1" = "27÷3=9, 0";
    "1,2" = "58÷2=29, 0";
    "1,3" = "22÷3=7, 1";
    "1,4" = "50÷6=8, 2";
    "1,5" = "16÷4=4, 0";

    "5,1" = "74÷9=8, 2";
    "5,2" = "24÷6=4, 0";
    "5,3" = "66÷7=9, 3";
    "5,4" = "19÷9=2, 1";
    "5,5" = "39÷3=13, 0";

    "9,1" = "24÷3=8, 0";
    "9,2" = "55÷5=11, 0";
    "9,3" = "77÷9=8, 5";
    "9,4" = "48÷5=9, 3";
    "9,5" = "99÷2=49, 1";

    "13,1" = "24÷8=3, 0";
    "13,2" = "26÷7=3, 5";
    "13,3" = "96÷4=24, 0";
    "13,4" = "62÷4=15, 2";
    "13,5" = "95÷8=11, 7";

    "17,1" = "29÷4=7, 1";
    "17,2" = "68÷8=8, 4";
    "17,3" = "85÷5=17, 0";
    "17,4" = "11÷4=2, 3";
    "17,5" = "79÷4=19, 3";
}

foreach ($key in $newValues.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $t.Cell($row, $col).Range.Text = $newValues[$key]
}
